$d = $word.ActiveDocument

# Find the range for "& Feiko Wielsma" within the subtitle paragraph
$rng = $d.Content
$found = $rng.Find.Execute("& Feiko Wielsma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng.Font.StrikeThrough = 1
